# Models Used.xlsx - add a new "K-Means Clustering / Credit Card Clustering" row
# right after the existing "Agglomerative (Hierarchical) Clustering / Relay States" row
# (new row 23; everything from the old row 23 onward shifts down by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Insert a new row at 23 - this shifts rows/merges/row-heights down by one,
#    which gets most of rows 23-28 (now 24-29) into the right place already.
# ---------------------------------------------------------------------------
$ws.Rows(23).Insert()

# ---------------------------------------------------------------------------
# 2. The "Unsupervised Learning - Clustering" merged label (old B23:B24) needs
#    to now span B23:B25. After the insert its text ended up on B24 (format
#    28 = merge-top), so move it up to B23 and make B24 a blank merge-middle
#    cell with the new border/wrap format, then restore the B25 merge-bottom
#    cell's format (it was left blank & correct after the insert).
# ---------------------------------------------------------------------------
$clusterLabel = $ws.Range("B24").Value2

$ws.Range("B24:B25").UnMerge()

$ws.Range("B24").Copy() | Out-Null
$ws.Range("B23").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("B23").Value2 = $clusterLabel

$ws.Range("B24").ClearContents()
$ws.Range("B27").Copy() | Out-Null                # format w/ borderId 10, center/center
$ws.Range("B24").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("B24").WrapText = $true

$ws.Range("B23:B25").Merge()

# ---------------------------------------------------------------------------
# 3. Fill in the new row's A/C cells (copy formatting from the row that
#    already carries the right style after the shift, then set the values).
# ---------------------------------------------------------------------------
$ws.Range("A24").Copy() | Out-Null
$ws.Range("A23").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false
$ws.Range("A23").Value2 = "K-Means Clustering"

$ws.Range("C24").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false
$ws.Range("C23").Value2 = "Credit Card Clustering"

# Row 24 (the old row 23, "Agglomerative (Hierarchical) Clustering") picked up
# an explicit row height in the authored edit.
$ws.Rows(24).RowHeight = 15.6

# ---------------------------------------------------------------------------
# 4. Hyperlinks: row insertion does not slide hyperlink ranges along with the
#    rows, so rebuild the whole collection - delete everything and re-add
#    each link against its (possibly shifted) cell, in original order, plus
#    the brand new one for C23.
# ---------------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C3"), "c. Jupyter Notebooks\Bike Rental Demand.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "c. Jupyter Notebooks\Temperature Trends.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C11"), "c. Jupyter Notebooks\Graduate Admissions.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C15"), "c. Jupyter Notebooks\Employee Attrition.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C16"), "c. Jupyter Notebooks\Cardiac Risk.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C18"), "c. Jupyter Notebooks\Bank Churn.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C19"), "c. Jupyter Notebooks\Skin Analysis.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C26"), "c. Jupyter Notebooks\Air Passenger.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C27"), "c. Jupyter Notebooks\Champagne Sales.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "c. Jupyter Notebooks\King County House Sales.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C13"), "c. Jupyter Notebooks\Bank Client Term Deposit.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C25"), "c. Jupyter Notebooks\National Health and Nutrition Examination Survey (NHANES).ipynb", "", "", "NHANES") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C29"), "c. Jupyter Notebooks\Oil Prices.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C6"), "c. Jupyter Notebooks\Laptop Price Analysis.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C28"), "c. Jupyter Notebooks\Mindtree Stock Price.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C7"), "c. Jupyter Notebooks\Possum Morphometric Analysis.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C8"), "c. Jupyter Notebooks\Possum Morphometric Analysis.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C17"), "c. Jupyter Notebooks\Breast Cancer.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C22"), "c. Jupyter Notebooks\Relay States.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C24"), "c. Jupyter Notebooks\Relay States.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C21"), "c. Jupyter Notebooks\Industrial System Monitoring.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C20"), "c. Jupyter Notebooks\Brain Stroke.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C12"), "c. Jupyter Notebooks\Loan Prediction.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C9"), "c. Jupyter Notebooks\Churn Prediction.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C10"), "c. Jupyter Notebooks\Stroke Prediction.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C14"), "c. Jupyter Notebooks\Chronic Kidney Disease Classification.ipynb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C23"), "c. Jupyter Notebooks\Credit Card Clustering.ipynb") | Out-Null

# ---------------------------------------------------------------------------
# 5. View state: selection moves to C23 (the new row).
# ---------------------------------------------------------------------------
$ws.Range("C23").Select()

Write-Output "done"
